$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Fill in the "Cost" (C) column values that were previously blank,
#    and fix a couple of "Quantity" (D) values. The "Total Cost" (E)
#    column recalculates automatically since it holds =C*D formulas.
# ---------------------------------------------------------------
$ws.Range("C2").Value = 70
$ws.Range("C3").Value = 100
$ws.Range("C4").Value = 55
$ws.Range("C5").Value = 45
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 0.3
$ws.Range("C7").Value = 0.75
$ws.Range("C8").Value = 10
$ws.Range("C9").Value = 1
$ws.Range("C10").Value = 5
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 2
$ws.Range("C12").Value = 1
$ws.Range("C13").Value = 5
$ws.Range("C16").Value = 2
$ws.Range("C17").Value = 60
$ws.Range("C18").Value = 0.75
$ws.Range("C19").Value = 0.5
$ws.Range("C20").Value = 2.5
$ws.Range("C21").Value = 25

# ---------------------------------------------------------------
# 2. Rename the connector part from "JST 8-pin" to
#    "GX16 - 8Pole connector" (row 17, column B).
# ---------------------------------------------------------------
$ws.Range("B17").Value = "GX16 – 8Pole connector"

# ---------------------------------------------------------------
# 3. Widen column B a little.
# ---------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 27.6

# ---------------------------------------------------------------
# 4. Insert two new rows right before the old TOTAL row (row 24),
#    pushing the old TOTAL row (and the blank row after it) down
#    by two rows.
# ---------------------------------------------------------------
$ws.Rows.Item(24).Insert()
$ws.Rows.Item(24).Insert()

# Re-use the formatting of the (old, now shifted) blank spacer row 23
# for the two newly inserted rows 24 and 25, copying cell-by-cell so
# every column keeps its own border/number-format combination.
foreach ($col in @("A","B","C","D","E","F","G","H")) {
  $ws.Range($col + "23").Copy()
  $ws.Range($col + "24").PasteSpecial(-4122)
  $ws.Range($col + "23").Copy()
  $ws.Range($col + "25").PasteSpecial(-4122)
}

# Row heights to match the rest of the "footer" block.
$ws.Rows.Item(24).RowHeight = 16.5
$ws.Rows.Item(25).RowHeight = 16.5
$ws.Rows.Item(27).RowHeight = 15.75

# ---------------------------------------------------------------
# 5. New row 24: "ALREADY BOUGHT"-style highlighted summary row for
#    items that still need to be designed/manufactured/tested.
# ---------------------------------------------------------------
$ws.Range("A24:E24").Interior.Color = 65535
$ws.Range("A24").Value = "<>"
$ws.Range("B24").Value = "Design, Manufacture and Testing"
$ws.Range("E24").Value = 4000

# Row 25 is left blank (same styling as the old spacer row 23).

# ---------------------------------------------------------------
# 6. Row 26 is the old TOTAL row, now shifted down two rows: update
#    its SUM formula so it covers the new data range.
# ---------------------------------------------------------------
$ws.Range("E26").Formula = "=SUM(E2:E25)"

# ---------------------------------------------------------------
# 7. Cosmetic: restore the active cell/selection like the saved file.
# ---------------------------------------------------------------
$ws.Range("D30").Select()
